$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2022")

# ---- Row 3 ----
# Columns are populated left to right so new shared-string entries are
# minted in the same order Excel would create them (C, D, E, F, G).
$ws.Cells.Item(3, 1).Value = 3
$ws.Cells.Item(3, 2).Value = 1
$ws.Cells.Item(3, 3).Value = "ABC"
$ws.Cells.Item(3, 4).Value = "13.07.2022"

# Column E ("12.08.2022") would otherwise be auto-recognised as a date
# serial by the engine's smart entry parser. Force it to stay text by
# temporarily marking the cell as Text before assigning the value, then
# restore the surrounding cell's original (bordered, General) look by
# copying the finished formatting from the neighbouring D column cell,
# which already carries the exact target style.
$ws.Cells.Item(3, 5).NumberFormat = "@"
$ws.Cells.Item(3, 5).Value = "12.08.2022"
$ws.Cells.Item(3, 4).Copy()
$ws.Cells.Item(3, 5).PasteSpecial(-4122)

$ws.Cells.Item(3, 6).Value = "Created"
$ws.Cells.Item(3, 7).Value = "C\....."

# ---- Row 4 (brand new row) ----
$ws.Cells.Item(4, 1).Value = 4
$ws.Cells.Item(4, 2).Value = 2
$ws.Cells.Item(4, 3).Value = "ABC"
$ws.Cells.Item(4, 4).Value = "13.07.2022"

$ws.Cells.Item(4, 5).NumberFormat = "@"
$ws.Cells.Item(4, 5).Value = "12.08.2022"
$ws.Cells.Item(4, 4).Copy()
$ws.Cells.Item(4, 5).PasteSpecial(-4122)

$ws.Cells.Item(4, 6).Value = "Created"
$ws.Cells.Item(4, 7).Value = "C\....."
